# Updates cryptocurrency price/volume figures on the sheet to match the
# latest scrape. All Price (D) and Volume(1h) (E) columns hold text values
# (not numbers) in the source data, e.g. "26.258.62" or "  -0.35%  ".
# For Price values that look like plain numbers (e.g. "212.86", "1.00"),
# a leading apostrophe is used (same as typing '212.86 directly into Excel)
# to force them to stay as text instead of being auto-converted to numeric
# values, matching the original inline-string formatting in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.258.62'
$ws.Range("D3").Value = '1.593.65'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Formula = '''212.86'
$ws.Range("D6").Formula = '''0.504'
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("D10").Formula = '''18.96'
$ws.Range("E10").Value = '  -2.55%  '
$ws.Range("D11").Formula = '''0.0851'
$ws.Range("E11").Value = '  +0.55%  '
$ws.Range("D12").Value = '1.817.84'
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").Value = '1.594.09'
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("E14").Value = '  -1.24%  '
$ws.Range("E15").Value = '  -2.53%  '
$ws.Range("D16").Formula = '''63.97'
$ws.Range("E16").Value = '  -1.11%  '
$ws.Range("D17").Value = '26.262.60'
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").Formula = '''7.37'
$ws.Range("E19").Value = '  -1.63%  '
$ws.Range("D20").Formula = '''214.17'
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("E24").Value = '  -3.06%  '
$ws.Range("D25").Formula = '''145.02'
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  -1.37%  '
$ws.Range("E28").Value = '  -0.38%  '
$ws.Range("E29").Value = '  -0.50%  '
$ws.Range("E30").Value = '  -2.48%  '
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("E32").Value = '  -0.48%  '
$ws.Range("D33").Value = '1.418.14'
$ws.Range("E33").Value = '  +5.70%  '
$ws.Range("D34").Formula = '''2.98'
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("E35").Value = '  -0.58%  '
$ws.Range("E36").Value = '  -1.10%  '
$ws.Range("E37").Value = '  -3.29%  '
$ws.Range("E38").Value = '  -0.78%  '
$ws.Range("E39").Value = '  +0.56%  '
$ws.Range("D40").Formula = '''5.80'
$ws.Range("E40").Value = '  +0.93%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Formula = '''0.967'
$ws.Range("E42").Value = '  -9.77%  '
$ws.Range("E43").Value = '  +0.95%  '
$ws.Range("D44").Formula = '''0.764'
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").Value = '1.730.07'
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Formula = '''60.92'
$ws.Range("E46").Value = '  -1.16%  '
$ws.Range("D47").Formula = '''86.98'
$ws.Range("E47").Value = '  -1.20%  '
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("E49").Value = '  -0.69%  '
$ws.Range("D50").Formula = '''0.0955'
$ws.Range("E50").Value = '  -3.08%  '
$ws.Range("D51").Formula = '''1.00'
$ws.Range("E51").Value = '  -0.05%  '
